# "Extended the microservice for course" - add two more professor rows
# (Gavrilut Dragos / Android, Lenuta Alboaie / Retele de calculatoare, Cloud
# Computing) to the Studenti sheet, each with a mailto hyperlink on the
# e-mail cell, matching the existing Florin Olariu row's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - fill in the same order the source workbook's shared-string table
# was built in (email, course, name, then the already-known "Prof." title).
$ws.Range("C3").Value = "dragos_gavrilut@gmail.com"
$ws.Range("D3").Value = "Android"
$ws.Range("A3").Value = "Gavriluț Dragoș"
$ws.Range("B3").Value = "Prof."

# Row 4
$ws.Range("A4").Value = "Lenuta Alboaie"
$ws.Range("B4").Value = "Lect. Dr."
$ws.Range("C4").Value = "lenuta_alboaie@gmail.com"
$ws.Range("D4").Value = "Retele de calculatoare, Cloud Computing"

# Hyperlink the two new e-mail cells the same way C2 already is, then
# re-apply the built-in "Hyperlink" cell style so it reuses the workbook's
# existing style entry instead of minting a new one.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:dragos_gavrilut@gmail.com")
$ws.Range("C3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:lenuta_alboaie@gmail.com")
$ws.Range("C4").Style = "Hyperlink"

# Widen the Email/Title->Email/Courses columns to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 37.8
$ws.Columns.Item(4).ColumnWidth = 36.5

# Move the active selection the way the authored workbook left it.
$ws.Range("D7").Select() | Out-Null
